# Apply numeric "want-to-go" count updates and a couple of text/URL
# corrections across all four worksheets, per the upstream data refresh.
$wb = $excel.ActiveWorkbook

# --- 展览 ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 251
$wsExhibit.Range("F7").Value = 12956
$wsExhibit.Range("F8").Value = 50
$wsExhibit.Range("F9").Value = 118
$wsExhibit.Range("F10").Value = 243
$wsExhibit.Range("F11").Value = 3014
$wsExhibit.Range("F13").Value = 6396
$wsExhibit.Range("F16").Value = 3378
$wsExhibit.Range("F18").Value = 161
$wsExhibit.Range("F20").Value = 37
$wsExhibit.Range("F23").Value = 30
$wsExhibit.Range("F24").Value = 3582
$wsExhibit.Range("F25").Value = 85
$wsExhibit.Range("F27").Value = 2748
$wsExhibit.Range("F28").Value = 2748
$wsExhibit.Range("F29").Value = 402
$wsExhibit.Range("F30").Value = 1877
$wsExhibit.Range("F32").Value = 210
$wsExhibit.Range("F33").Value = 6576
$wsExhibit.Range("F35").Value = 165
$wsExhibit.Range("F36").Value = 739
$wsExhibit.Range("F37").Value = 1975
$wsExhibit.Range("F38").Value = 1292
$wsExhibit.Range("F39").Value = 95
$wsExhibit.Range("F40").Value = 1032
$wsExhibit.Range("F43").Value = 219
$wsExhibit.Range("F44").Value = 1146
$wsExhibit.Range("F46").Value = 127
$wsExhibit.Range("F47").Value = 1193
$wsExhibit.Range("F48").Value = 1760
$wsExhibit.Range("F49").Value = 156
$wsExhibit.Range("F50").Value = 1168

# --- 演出 ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F15").Value = 96

# --- 本地生活 ---
$wsLocalLife = $wb.Worksheets.Item("本地生活")
$wsLocalLife.Range("F2").Value = 425
$wsLocalLife.Range("F3").Value = 589

# --- 全部类型 ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 425
$wsAll.Range("F7").Value = 589
$wsAll.Range("F8").Value = 251
$wsAll.Range("F10").Value = 12956
$wsAll.Range("F11").Value = 118
$wsAll.Range("F13").Value = 243
$wsAll.Range("F14").Value = 3014
$wsAll.Range("F15").Value = 6397
$wsAll.Range("F17").Value = 3377
$wsAll.Range("F19").Value = 161
$wsAll.Range("F21").Value = 37
$wsAll.Range("F25").Value = 30
$wsAll.Range("F26").Value = 3582
$wsAll.Range("F28").Value = 2748
$wsAll.Range("F29").Value = 402
$wsAll.Range("F30").Value = 1877
$wsAll.Range("F32").Value = 210
$wsAll.Range("F33").Value = 6576
$wsAll.Range("F34").Value = 96
$wsAll.Range("F36").Value = 165
$wsAll.Range("F37").Value = 739
$wsAll.Range("F38").Value = 1975
$wsAll.Range("F40").Value = 1292
$wsAll.Range("F41").Value = 95
$wsAll.Range("F42").Value = 1032
$wsAll.Range("F44").Value = 219
$wsAll.Range("F45").Value = 1146
$wsAll.Range("F46").Value = 127
$wsAll.Range("F48").Value = 1760
$wsAll.Range("F50").Value = 156

# --- Text/URL corrections (also appear on 展览 and 全部类型) ---
$wsExhibit.Range("C32").Value = "北京·原神x星穹铁道x绝区零同人only"
$wsExhibit.Range("I31").Value = "//i2.hdslb.com/bfs/openplatform/202407/tybJZC5s1721724558989.jpeg"
$wsAll.Range("C32").Value = "北京·原神x星穹铁道x绝区零同人only"
$wsAll.Range("I31").Value = "//i2.hdslb.com/bfs/openplatform/202407/tybJZC5s1721724558989.jpeg"
